$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-13 14:21:54"
$wsZhCn.Range("H2").Value = "2016-03-13 14:22:12"
$wsZhCn.Range("E5").Value = "2016-03-13 14:21:54"
$wsZhCn.Range("H5").Value = "2016-03-13 14:22:12"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-13 14:21:58"
$wsDeDe.Range("H2").Value = "2016-03-13 14:22:19"
$wsDeDe.Range("E5").Value = "2016-03-13 14:21:58"
$wsDeDe.Range("H5").Value = "2016-03-13 14:22:19"
